$p = $ppt.ActivePresentation

# --- Slide 2: "ZoneTexte 5" -------------------------------------------------
# Add a new bullet paragraph after "... de la grille dans la direction donnée."
# and resize the text box (spAutoFit grows the shape to fit the extra line).
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(4)
$tr2 = $sh2.TextFrame.TextRange
$para3 = $tr2.Paragraphs(3, 1)
$para3.InsertAfter([char]13 + "Si deux tuiles adjacentes sont identiques dans la même direction, elles s’additionnent. ") | Out-Null

$sh2.Left = 510
$sh2.Top = 114.37543307086614
$sh2.Width = 428.25
$sh2.Height = 334.4344094488189

# --- Slide 4: "ZoneTexte 4" -------------------------------------------------
# Fix the typo in the rotateMatrix description.
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(4)
$tr4 = $sh4.TextFrame.TextRange
$para2 = $tr4.Paragraphs(2)
$para2.Runs(1).Text = "fait tourner de 90° dans sens trigo la matrice"
